$d = $word.ActiveDocument

# The footer block at the end of the document ("Ver no Jupiter Salvar em pdf
# Salvar em docx" + the copyright line), along with the blank paragraph that
# precedes it, is being removed. Locate those two text paragraphs by their
# content and delete the range spanning from the start of the blank
# paragraph immediately before the "Ver no Jupiter" paragraph through the
# end of the copyright paragraph (inclusive of its paragraph mark).

$paras = $d.Paragraphs
$count = $paras.Count

$jupiterIdx = $null
$copyrightIdx = $null

for ($i = 1; $i -le $count; $i++) {
    $txt = $paras.Item($i).Range.Text
    if ($jupiterIdx -eq $null -and $txt -like "*Ver no Jupiter*") {
        $jupiterIdx = $i
    }
    if ($copyrightIdx -eq $null -and $txt -like "*Contact: luizeleno@usp.br*") {
        $copyrightIdx = $i
    }
}

if ($jupiterIdx -ne $null -and $copyrightIdx -ne $null) {
    # The empty paragraph that immediately precedes the "Ver no Jupiter" one
    # is also removed as part of this edit.
    $firstIdx = $jupiterIdx - 1

    $startRange = $paras.Item($firstIdx).Range.Start
    $endRange = $paras.Item($copyrightIdx).Range.End

    $d.Range($startRange, $endRange).Delete()
}
